# Ran player stats and added the newly-compatible players
# (Alexander Mattison, Austin Ekeler) to the RB aggregate sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reuse the existing row formatting (green fill for the first player block,
# yellow fill for the second) instead of creating brand-new styles, by
# copying formats from the existing rows down onto the new ones.
$ws.Range("A2:F4").Copy()
$ws.Range("A8:F10").PasteSpecial(-4122)

$ws.Range("A5:F7").Copy()
$ws.Range("A11:F13").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Alexander Mattison
$ws.Range("A8").Value = "Alexander Mattison"
$ws.Range("B8").Value = "Group1"
$ws.Range("C8").Value = 4.266666666666667
$ws.Range("D8").Value = 8.299999999999999
$ws.Range("E8").Value = 41.93333333333333
$ws.Range("F8").Value = 23.66666666666667

$ws.Range("A9").Value = "Alexander Mattison"
$ws.Range("B9").Value = "Group2"
$ws.Range("C9").Value = 3.633333333333333
$ws.Range("D9").Value = 6.899999999999999
$ws.Range("E9").Value = 42.63333333333333
$ws.Range("F9").Value = 24.66666666666667

$ws.Range("A10").Value = "Alexander Mattison"
$ws.Range("B10").Value = "Difference"
$ws.Range("C10").Value = -0.6333333333333337
$ws.Range("D10").Value = -1.399999999999999
$ws.Range("E10").Value = 0.7000000000000028
$ws.Range("F10").Value = 1

# Austin Ekeler
$ws.Range("A11").Value = "Austin Ekeler"
$ws.Range("B11").Value = "Group1"
$ws.Range("C11").Value = 4.4
$ws.Range("D11").Value = 9.166666666666666
$ws.Range("E11").Value = 47.86666666666667
$ws.Range("F11").Value = 37

$ws.Range("A12").Value = "Austin Ekeler"
$ws.Range("B12").Value = "Group2"
$ws.Range("C12").Value = 4.266666666666667
$ws.Range("D12").Value = 8.566666666666666
$ws.Range("E12").Value = 49.13333333333333
$ws.Range("F12").Value = 33

$ws.Range("A13").Value = "Austin Ekeler"
$ws.Range("B13").Value = "Difference"
$ws.Range("C13").Value = -0.1333333333333337
$ws.Range("D13").Value = -0.5999999999999996
$ws.Range("E13").Value = 1.266666666666659
$ws.Range("F13").Value = -4
